# Add a new "2022-Q1" sheet (fund-holding detail) before the "总计" (total)
# sheet, and add a corresponding summary row to "总计".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the "2022-Q1" worksheet by duplicating "2021-Q4" (same header
#    row / column layout / cell styles) and dropping it right before 总计.
# ---------------------------------------------------------------------
$templateSheet = $wb.Worksheets.Item("2021-Q4")
$totalSheet    = $wb.Worksheets.Item("总计")
$templateSheet.Copy($totalSheet)

# The freshly made copy is selected/active right after Copy(); grab it by
# the name Excel gives duplicates ("<source> (2)").
$q1 = $wb.Worksheets.Item("2021-Q4 (2)")
$q1.Name = "2022-Q1"

# ---------------------------------------------------------------------
# 2. Fill in the fund-holding rows for 2022-Q1.
# ---------------------------------------------------------------------
# columns: A idx(n) | B code(text) | C name(text) | D size(text) |
#          E stock position(text) | F position ratio(text) |
#          G market value(text) | H rank(n)
$rows = @(
    @(0,  "516950", "银华中证基建交易型开放式指数证券投资基金", "10.41", "97.55", "4.22", "0.4393", 7),
    @(1,  "320011", "诺安中小盘精选混合",                         "3.67",  "84.64", "5.22", "0.1916", 2),
    @(2,  "001528", "诺安先进制造股票",                           "2.44",  "87.44", "5.61", "0.1369", 3),
    @(3,  "006977", "农银汇理海棠三年定期开放混合",               "4.64",  "63.10", "2.45", "0.1137", 9),
    @(4,  "320015", "诺安行业轮动混合",                           "1.29",  "85.68", "4.53", "0.0584", 4),
    @(5,  "006429", "诺安恒鑫混合",                               "0.82",  "85.51", "4.78", "0.0392", 5),
    @(6,  "000646", "华润元大量化优选混合A",                      "0.80",  "67.15", "3.17", "0.0254", 10),
    @(7,  "515870", "嘉实中证先进制造100策略ETF",                 "0.42",  "98.79", "4.40", "0.0185", 8),
    @(8,  "007827", "华润元大量化优选混合C",                      "0.33",  "67.15", "3.17", "0.0105", 10),
    @(9,  "090011", "大成核心双动力混合",                         "0.34",  "93.14", "2.05", "0.0070", 7),
    @(10, "161718", "招商沪深300高贝塔指数",                      "0.20",  "94.52", "1.36", "0.0027", 7),
    @(11, "519165", "新华鑫利灵活配置混合",                       "0.05",  "74.84", "4.01", "0.0020", 1)
)

$lastRow = 1 + $rows.Count   # header is row 1

# The copy of 2021-Q4 only had 6 data rows (rows 2-7); extend the
# column-A style (bold/centered/bordered) down to every new data row so
# rows 8..13 match rows 2..7 before any values are written.
$q1.Range("A7").Copy()
$q1.Range("A8:A$lastRow").PasteSpecial(-4122)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $data = $rows[$i]

    $q1.Cells.Item($r, 1).Value = $data[0]

    # Fund code / size / position / ratio / value look numeric, so force
    # Text format before assignment or Excel will silently convert them
    # to numbers (and mangle leading zeros in fund codes).
    $q1.Cells.Item($r, 2).NumberFormat = "@"
    $q1.Cells.Item($r, 2).Value = $data[1]

    $q1.Cells.Item($r, 3).Value = $data[2]

    $q1.Cells.Item($r, 4).NumberFormat = "@"
    $q1.Cells.Item($r, 4).Value = $data[3]

    $q1.Cells.Item($r, 5).NumberFormat = "@"
    $q1.Cells.Item($r, 5).Value = $data[4]

    $q1.Cells.Item($r, 6).NumberFormat = "@"
    $q1.Cells.Item($r, 6).Value = $data[5]

    $q1.Cells.Item($r, 7).NumberFormat = "@"
    $q1.Cells.Item($r, 7).Value = $data[6]

    $q1.Cells.Item($r, 8).Value = $data[7]
}

# ---------------------------------------------------------------------
# 3. Add the 2022-Q1 summary row to "总计" (insert as the new row 2, push
#    the rest down, and renumber the leading index column).
# ---------------------------------------------------------------------
# Re-resolve "总计" by name: the sheet collection shifted when the new
# "2022-Q1" sheet was inserted in front of it, so the old $totalSheet
# handle now refers to the wrong (positionally shifted) worksheet.
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

# Clear the formatting Excel auto-applies to an inserted row's cells
# (it copies the row above) for the text/number columns, then restore
# the bold/centered/bordered style on column A to match the other rows.
$totalSheet.Range("B2:D2").ClearFormats()
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 12
$totalSheet.Cells.Item(2, 4).Value = 1.05

# Renumber the old rows' index column (was 0..4, now 1..5).
$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(4, 1).Value = 2
$totalSheet.Cells.Item(5, 1).Value = 3
$totalSheet.Cells.Item(6, 1).Value = 4
$totalSheet.Cells.Item(7, 1).Value = 5

Write-Output "done"
